$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 123
$ws.Range("A3").Value = 123
$ws.Range("A4").Value = 123
